$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for the new authorized person
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "MEJIA ARANGO"
$ws.Range("C3").Value = "ISABELLA MARIA"
$ws.Range("D3").Value = "T1019906212"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Activo(a)"
$ws.Range("G3").Value = 1019906212
$ws.Range("H3").Value = "F Mejia"
$ws.Range("I3").Value = "1019906212.jpg"

# Fix the ID_STATUS and STATUS for existing row (now the person is active)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "Activo(a)"

# Adjust column width for column I (imageUrl) to match new content
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666

# Update selection to reflect F2 as the active cell
$ws.Range("F2").Select()
